$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 formula change
$ws.Range("D5").Formula = "=(C5+C4)/4799"

# Row 6 formula change
$ws.Range("D6").Formula = "=(C6+C5+C4)/4799"

# Row 7 values and formula
$ws.Range("B7").Value = 4799
$ws.Range("C7").Value = 481
$ws.Range("D7").Formula = "=(C7+C6+C5+C4)/B14"
$ws.Range("F7").Value = 61
$ws.Range("G7").Value = 298
$ws.Range("H7").Value = 720
$ws.Range("I7").Value = 1965

# Row 8 - C8 gets a single space string
$ws.Range("C8").Value = " "

# Row 14
$ws.Range("B14").Value = 4799
$ws.Range("C14").Formula = "=SUM(C4:C13)"
$ws.Range("D14").Value = 0.58762242133777876
$ws.Range("D14").NumberFormat = "0.00%"

# Row 15 - new row, D15 gets single space string
$ws.Range("D15").Value = " "

# Selection change
$ws.Range("A1:J16").Select()
